$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 5 (P0_2) is being re-themed to the "external" (blue) row coloring used by
# rows 3/4/11... - copy that formatting over, which leaves the existing
# cell values/content untouched.
$ws.Range("B3:F3").Copy()
$ws.Range("B5:F5").PasteSpecial(-4122)

# Row 7 (P0_4) moves from the "on-board" (orange) coloring to the "input"
# (green) coloring used by rows 6/8/9 - copy that formatting over too.
$ws.Range("B6:F6").Copy()
$ws.Range("B7:F7").PasteSpecial(-4122)

# Row 5 now has a previously unused GPIO - OLED CS.
$ws.Range("C5").Value = "OLED CS"

# The relay's GPIO (D7, "LED Green") moves off this pin; add a clarifying
# note instead (4 cluster for relay => z2m).
$ws.Range("D7").Value = ""
$ws.Range("G7").Value = "LED Green?"

# Move the active selection to the relay's GPIO cell.
$ws.Range("D11").Select()
